$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DEC")

$ws.Range("A3").Value = "feb/23"
$ws.Range("A5").Value = "apr/23"
$ws.Range("A6").Value = "may/23"
$ws.Range("A9").Value = "aug/23"
$ws.Range("A10").Value = "sep/23"
$ws.Range("A11").Value = "oct/23"
$ws.Range("A13").Value = "dec/23"
$ws.Range("A15").Value = "feb/24"
$ws.Range("A17").Value = "apr/24"
$ws.Range("A18").Value = "may/24"
$ws.Range("A21").Value = "aug/24"

$ws.Range("E8:E9").Select()
